$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.75
$ws.Range("J2").Value = 3.2
$ws.Range("L2").Value = 4.75
$ws.Range("M2").Value = 1.1
$ws.Range("O2").Value = 1.58
$ws.Range("Q2").Value = 3.1
$ws.Range("R2").Value = 1.36
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 9
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 23
$ws.Range("AF2").Value = 101
$ws.Range("AI2").Value = 17
$ws.Range("AJ2").Value = 15
$ws.Range("AW2").Value = 5.5
$ws.Range("AX2").Value = 23
$ws.Range("BB2").Value = 501

# Row 3
$ws.Range("A3").Value = "0b9yyKnJ"
$ws.Range("C3").Value = "23:30"
$ws.Range("D3").Value = "AUSTRALIA - A-LEAGUE"
$ws.Range("E3").Value = "Wellington Phoenix"
$ws.Range("F3").Value = "Melbourne Victory"
$ws.Range("G3").Value = 4.5
$ws.Range("H3").Value = 3.75
$ws.Range("I3").Value = 1.8
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 2.1
$ws.Range("L3").Value = 2.4
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 9.5
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 3.4
$ws.Range("Q3").Value = 2.05
$ws.Range("R3").Value = 1.8
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 1.95
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 11
$ws.Range("X3").Value = 21
$ws.Range("Y3").Value = 15
$ws.Range("Z3").Value = 51
$ws.Range("AA3").Value = 41
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 9.5
$ws.Range("AD3").Value = 7
$ws.Range("AH3").Value = 6.5
$ws.Range("AI3").Value = 8
$ws.Range("AJ3").Value = 8.5
$ws.Range("AK3").Value = 13
$ws.Range("AL3").Value = 15
$ws.Range("AM3").Value = 29
$ws.Range("AN3").Value = 6.5
$ws.Range("AO3").Value = 26
$ws.Range("AP3").Value = 34
$ws.Range("AQ3").Value = 81
$ws.Range("AR3").Value = 126
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.63
$ws.Range("AW3").Value = 3.75
$ws.Range("AX3").Value = 9.5
$ws.Range("AY3").Value = 21
$ws.Range("AZ3").Value = 34
$ws.Range("BA3").Value = 51
$ws.Range("BB3").Value = 151
$ws.Range("BC3").Value = 501
$ws.Range("BD3").Value = 126

# Row 5
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 11
